# Apply the "Scenario 10" edits:
#  - Decrease probability of fetal death before 4 weeks from conception
#    (Phase1!B2:B5 and Phase2!C2:C5), which cascades into the dependent
#    "no loss" formulas in column D (Phase1) / E (Phase2).
#  - Update the saved selection/active-sheet state to reflect the cells
#    that were last touched on each sheet.

$wb = $excel.ActiveWorkbook

$phase1 = $wb.Worksheets.Item("Phase1")
$phase2 = $wb.Worksheets.Item("Phase2")

# --- Update probabilities on Phase1 ---
$phase1.Range("B2").Value = 0.1
$phase1.Range("B3").Value = 0.1
$phase1.Range("B4").Value = 0.05
$phase1.Range("B5").Value = 0.05

# --- Update probabilities on Phase2 ---
$phase2.Range("C2").Value = 0.1
$phase2.Range("C3").Value = 0.1
$phase2.Range("C4").Value = 0.05
$phase2.Range("C5").Value = 0.05

# --- Replicate the editor's on-screen selection / active sheet state ---
# The user last selected C2:C5 on Phase2 ...
[void]$phase2.Activate()
[void]$phase2.Range("C2:C5").Select()

# ... then moved to Phase1 and selected B2:B5, leaving Phase1 active.
[void]$phase1.Activate()
[void]$phase1.Range("B2:B5").Select()
